$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.707.53'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.292.19'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.26'
$ws.Range('E5').Value = '  +18.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '268.82'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.613'
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.46'
$ws.Range('E10').Value = '  +4.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0936'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.52'
$ws.Range('E12').Value = '  +9.31%  '
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.59'
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.631.54'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.848'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.286.33'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.605.94'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.59'
$ws.Range('E20').Value = '  +6.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.36'
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('E22').Value = '  +2.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.20'
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.51'
$ws.Range('E24').Value = '  +4.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.84'
$ws.Range('E25').Value = '  +14.27%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.46'
$ws.Range('E27').Value = '  +2.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.33'
$ws.Range('E28').Value = '  +5.74%  '
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.39'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.63'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0923'
$ws.Range('E33').Value = '  +4.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.52'
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.71'
$ws.Range('E36').Value = '  +8.04%  '
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0354'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.79'
$ws.Range('E39').Value = '  +11.92%  '
$ws.Range('E40').Value = '  +4.34%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.242'
$ws.Range('E41').Value = '  +2.68%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.26'
$ws.Range('E42').Value = '  +14.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.69'
$ws.Range('E43').Value = '  +11.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.42'
$ws.Range('E44').Value = '  +5.13%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.94'
$ws.Range('E46').Value = '  +13.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.74'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.41'
$ws.Range('E48').Value = '  +5.31%  '
$ws.Range('E49').Value = '  -1.69%  '
$ws.Range('E50').Value = '  +3.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.443'
$ws.Range('E51').Value = '  +3.57%  '
